$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2
$ws.Range("R2").Value = 14
$ws.Range("S2").Value = 22
$ws.Range("T2").Value = 32
$ws.Range("U2").Value = 127
$ws.Range("V2").Value = 25.93354430379747
$ws.Range("W2").Value = 18.59421350768726
$ws.Range("Y2").Value = 24
$ws.Range("Z2").Value = 38
$ws.Range("AA2").Value = 53
$ws.Range("AB2").Value = 231
$ws.Range("AC2").Value = 47.5759493670886
$ws.Range("AD2").Value = 39.54339297110542
$ws.Range("AF2").Value = 2
$ws.Range("AI2").Value = 15
$ws.Range("AJ2").Value = 3.139240506329114
$ws.Range("AK2").Value = 2.773956785304159

# Row 5
$ws.Range("R5").Value = 11
$ws.Range("S5").Value = 19
$ws.Range("T5").Value = 29
$ws.Range("U5").Value = 90
$ws.Range("V5").Value = 21.94303797468354
$ws.Range("W5").Value = 15.23706578047602
$ws.Range("Y5").Value = 26
$ws.Range("Z5").Value = 37
$ws.Range("AA5").Value = 53.25
$ws.Range("AB5").Value = 174
$ws.Range("AC5").Value = 45.81645569620253
$ws.Range("AD5").Value = 32.42374959615592
$ws.Range("AF5").Value = 2
$ws.Range("AH5").Value = 3
$ws.Range("AI5").Value = 15
$ws.Range("AJ5").Value = 3.09493670886076
$ws.Range("AK5").Value = 2.567904400840213

# Row 6
$ws.Range("R6").Value = 3
$ws.Range("S6").Value = 7
$ws.Range("T6").Value = 12
$ws.Range("U6").Value = 65
$ws.Range("V6").Value = 9.79746835443038
$ws.Range("W6").Value = 10.04419408615073
$ws.Range("Y6").Value = 0
$ws.Range("Z6").Value = 8
$ws.Range("AA6").Value = 11
$ws.Range("AB6").Value = 43
$ws.Range("AC6").Value = 8.018987341772151
$ws.Range("AD6").Value = 9.24823210492178
$ws.Range("AF6").Value = 0
$ws.Range("AG6").Value = 1
$ws.Range("AH6").Value = 1
$ws.Range("AI6").Value = 5
$ws.Range("AJ6").Value = 0.8069620253164557
$ws.Range("AK6").Value = 0.9036516458049771

# Row 7
$ws.Range("R7").Value = 15
$ws.Range("S7").Value = 25
$ws.Range("T7").Value = 34
$ws.Range("U7").Value = 108
$ws.Range("V7").Value = 25.67219917012448
$ws.Range("W7").Value = 13.9754403231824
$ws.Range("Y7").Value = 26
$ws.Range("Z7").Value = 36
$ws.Range("AA7").Value = 49
$ws.Range("AB7").Value = 174
$ws.Range("AC7").Value = 39.25311203319502
$ws.Range("AD7").Value = 21.75350164973208
$ws.Range("AF7").Value = 2
$ws.Range("AG7").Value = 3
$ws.Range("AH7").Value = 3
$ws.Range("AI7").Value = 11
$ws.Range("AJ7").Value = 2.809128630705394
$ws.Range("AK7").Value = 1.608855593546295

# Row 10
$ws.Range("R10").Value = 15
$ws.Range("S10").Value = 23
$ws.Range("T10").Value = 31
$ws.Range("V10").Value = 23.28630705394191
$ws.Range("W10").Value = 10.94745876396246
$ws.Range("Y10").Value = 24
$ws.Range("Z10").Value = 35
$ws.Range("AA10").Value = 47
$ws.Range("AB10").Value = 141
$ws.Range("AC10").Value = 36.87966804979253
$ws.Range("AD10").Value = 20.00453181713292
$ws.Range("AF10").Value = 2
$ws.Range("AH10").Value = 3
$ws.Range("AI10").Value = 10
$ws.Range("AJ10").Value = 2.701244813278008
$ws.Range("AK10").Value = 1.597510597977387

# Row 11
$ws.Range("R11").Value = 6
$ws.Range("S11").Value = 12
$ws.Range("T11").Value = 19
$ws.Range("U11").Value = 112
$ws.Range("V11").Value = 15.01244813278008
$ws.Range("W11").Value = 14.51019220174818
$ws.Range("Y11").Value = 0
$ws.Range("Z11").Value = 9
$ws.Range("AA11").Value = 17
$ws.Range("AB11").Value = 35
$ws.Range("AC11").Value = 9.95850622406639
$ws.Range("AD11").Value = 8.95441815117
$ws.Range("AF11").Value = 0
$ws.Range("AG11").Value = 1
$ws.Range("AH11").Value = 2
$ws.Range("AI11").Value = 4
$ws.Range("AJ11").Value = 1.04149377593361
$ws.Range("AK11").Value = 0.9210163368091974
